$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.43180100408847
$ws.Range("C2").Value = 0.5561946450787332
$ws.Range("D2").Value = 0.01898470219389026
$ws.Range("E2").Value = 0.418868025178611
$ws.Range("F2").Value = 1.22420842577543
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.6851918156816694
$ws.Range("N2").Value = 0.9873592261007289
$ws.Range("B3").Value = 1.270843218843311
$ws.Range("C3").Value = 0.4868603933235249
$ws.Range("D3").Value = 0.01914739503100904
$ws.Range("E3").Value = 0.3649584804624908
$ws.Range("F3").Value = 1.168444706594499
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.6700495542047094
$ws.Range("N3").Value = 1.003888397040395
$ws.Range("B4").Value = 1.172501167612211
$ws.Range("C4").Value = 0.4444327630589555
$ws.Range("D4").Value = 0.01926975491615224
$ws.Range("E4").Value = 0.3320111203335046
$ws.Range("F4").Value = 1.135248235746545
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.6613475268064306
$ws.Range("N4").Value = 1.014555381492585
$ws.Range("B5").Value = 1.132544577260489
$ws.Range("C5").Value = 0.4271769619589918
$ws.Range("D5").Value = 0.0193251675234194
$ws.Range("E5").Value = 0.3186199087233206
$ws.Range("F5").Value = 1.121978765897921
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.6579493326846304
$ws.Range("N5").Value = 1.01903186804314
$ws.Range("B6").Value = 1.125916869622642
$ws.Range("C6").Value = 0.4243136222090698
$ws.Range("D6").Value = 0.01933470107326229
$ws.Range("E6").Value = 0.3163983286809184
$ws.Range("F6").Value = 1.119790873934249
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.6573939500089239
$ws.Range("N6").Value = 1.019782994033985
$ws.Range("B7").Value = 1.171961823827075
$ws.Range("C7").Value = 0.4441999115416593
$ws.Range("D7").Value = 0.01927047988786512
$ws.Range("E7").Value = 0.331830384154614
$ws.Range("F7").Value = 1.135068237959004
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.661301100826023
$ws.Range("N7").Value = 1.014615229045823
$ws.Range("B8").Value = 1.376199399961763
$ws.Range("C8").Value = 0.5322568711487179
$ws.Range("D8").Value = 0.01903607571363963
$ws.Range("E8").Value = 0.4002460886429873
$ws.Range("F8").Value = 1.204762350224854
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.6798461201412564
$ws.Range("N8").Value = 0.9929505784159041
$ws.Range("B9").Value = 1.780764153406892
$ws.Range("C9").Value = 0.7061947151043455
$ws.Range("D9").Value = 0.01875920105149831
$ws.Range("E9").Value = 0.535788821193421
$ws.Range("F9").Value = 1.349887909244501
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.7210197532924383
$ws.Range("N9").Value = 0.9546052659177207
$ws.Range("B10").Value = 2.08077230338705
$ws.Range("C10").Value = 0.8349296160035919
$ws.Range("D10").Value = 0.01867358378142114
$ws.Range("E10").Value = 0.636450542645548
$ws.Range("F10").Value = 1.461933247349464
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.7543200522087972
$ws.Range("N10").Value = 0.9289952987538861
$ws.Range("B11").Value = 2.217924869523586
$ws.Range("C11").Value = 0.8937386091240569
$ws.Range("D11").Value = 0.01866162578351194
$ws.Range("E11").Value = 0.6825312684266436
$ws.Range("F11").Value = 1.514142216851695
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.7701580401969039
$ws.Range("N11").Value = 0.9179096619487979
$ws.Range("B12").Value = 2.269963472863196
$ws.Range("C12").Value = 0.9160465645317117
$ws.Range("D12").Value = 0.0186610976509769
$ws.Range("E12").Value = 0.7000266628881207
$ws.Range("F12").Value = 1.534095255452371
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.7762567557253419
$ws.Range("N12").Value = 0.9137937574140693
$ws.Range("B13").Value = 2.2587514234379
$ws.Range("C13").Value = 0.9112404005356325
$ws.Range("D13").Value = 0.0186610315642568
$ws.Range("E13").Value = 0.6962566223281641
$ws.Range("F13").Value = 1.529789816759688
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.7749387535038181
$ws.Range("N13").Value = 0.9146765315895671
$ws.Range("B14").Value = 2.2222040482431
$ws.Range("C14").Value = 0.8955731150404063
$ws.Range("D14").Value = 0.01866150147474244
$ws.Range("E14").Value = 0.6839696872497427
$ws.Range("F14").Value = 1.515780079165381
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.7706577447913077
$ws.Range("N14").Value = 0.9175693961619196
$ws.Range("B15").Value = 2.199831132061433
$ws.Range("C15").Value = 0.8859815169456056
$ws.Range("D15").Value = 0.01866231381580974
$ws.Range("E15").Value = 0.6764496476126993
$ws.Range("F15").Value = 1.507222628932141
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.7680487461037444
$ws.Range("N15").Value = 0.9193520606761254
$ws.Range("B16").Value = 2.071823072738027
$ws.Range("C16").Value = 0.8310915125155702
$ws.Range("D16").Value = 0.01867491848107505
$ws.Range("E16").Value = 0.633445221417773
$ws.Range("F16").Value = 1.458546580191637
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.7532990526687229
$ws.Range("N16").Value = 0.929731199063383
$ws.Range("B17").Value = 1.993470925925237
$ws.Range("C17").Value = 0.7974834466871812
$ws.Range("D17").Value = 0.01868964761732528
$ws.Range("E17").Value = 0.6071403153444521
$ws.Range("F17").Value = 1.429005889067014
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.7444286503622948
$ws.Range("N17").Value = 0.9362436128909799
$ws.Range("B18").Value = 1.948468322837186
$ws.Range("C18").Value = 0.7781760405782165
$ws.Range("D18").Value = 0.01870065255525333
$ws.Range("E18").Value = 0.592037403652057
$ws.Range("F18").Value = 1.412131211511991
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.7393914045231753
$ws.Range("N18").Value = 0.9400424614798908
$ws.Range("B19").Value = 1.933241982149468
$ws.Range("C19").Value = 0.7716427625505276
$ws.Range("D19").Value = 0.01870481012150549
$ws.Range("E19").Value = 0.5869283286249072
$ws.Range("F19").Value = 1.40643758340525
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.737696942551608
$ws.Range("N19").Value = 0.9413377783360701
$ws.Range("B20").Value = 2.001805043607987
$ws.Range("C20").Value = 0.8010586757072247
$ws.Range("D20").Value = 0.0186878166532054
$ws.Range("E20").Value = 0.6099376967771946
$ws.Range("F20").Value = 1.432138472512435
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.7453662010517519
$ws.Range("N20").Value = 0.935544855845535
$ws.Range("B21").Value = 2.232936096178037
$ws.Range("C21").Value = 0.9001739156577742
$ws.Range("D21").Value = 0.01866125392796647
$ws.Range("E21").Value = 0.6875773871021522
$ws.Range("F21").Value = 1.519890085992898
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.7719124173112704
$ws.Range("N21").Value = 0.9167174600614842
$ws.Range("B22").Value = 2.38458974426311
$ws.Range("C22").Value = 0.9651758285419874
$ws.Range("D22").Value = 0.01866726640396266
$ws.Range("E22").Value = 0.7385872800035287
$ws.Range("F22").Value = 1.578307509139734
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.7898528354055259
$ws.Range("N22").Value = 0.9048909335181037
$ws.Range("B23").Value = 2.303593272698947
$ws.Range("C23").Value = 0.9304616225160771
$ws.Range("D23").Value = 0.01866187882006187
$ws.Range("E23").Value = 0.7113364733707357
$ws.Range("F23").Value = 1.547029931490101
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.7802229355520041
$ws.Range("N23").Value = 0.9111589380551806
$ws.Range("B24").Value = 1.998037055126417
$ws.Range("C24").Value = 0.7994422702759607
$ws.Range("D24").Value = 0.01868863653876218
$ws.Range("E24").Value = 0.6086729382376603
$ws.Range("F24").Value = 1.430721893424504
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.7449421399710729
$ws.Range("N24").Value = 0.9358605933551623
$ws.Range("B25").Value = 1.670853665900438
$ws.Range("C25").Value = 0.6589872455524528
$ws.Range("D25").Value = 0.01881388412277829
$ws.Range("E25").Value = 0.4989489557592321
$ws.Range("F25").Value = 1.309695057901877
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.7093543365496231
$ws.Range("N25").Value = 0.9645315964188583
